$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the content of the last publication row (A21:C21), and clear the
# title text in D21 while keeping its italic style, effectively dropping
# the "Family Dynamics, Birth Timing, and Child Temperament..." entry.
$ws.Range("A21:D21").ClearContents()

# Update the active selection to D22 (the now-empty row after the last entry)
$ws.Range("D22").Select()
